$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so values like "1.00" or
# multi-dot numbers (e.g. "68.724.73") are stored verbatim, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '68.724.73'
$ws.Range("E2").Value = '  +1.83%  '

$ws.Range("D3").Value = '3.784.46'
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '596.87'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").Value = '169.61'
$ws.Range("E6").Value = '  -0.40%  '

$ws.Range("D7").Value = '3.783.75'
$ws.Range("E7").Value = '  +0.62%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  -1.69%  '

$ws.Range("D11").Value = '6.53'
$ws.Range("E11").Value = '  +0.60%  '

$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("D13").Value = '0.0000264'
$ws.Range("E13").Value = '  -3.08%  '

$ws.Range("D14").Value = '36.92'
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("D15").Value = '4.420.12'
$ws.Range("E15").Value = '  +0.64%  '

$ws.Range("D16").Value = '3.785.02'
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").Value = '68.686.00'
$ws.Range("E17").Value = '  +1.58%  '

$ws.Range("D18").Value = '18.18'
$ws.Range("E18").Value = '  -3.66%  '

$ws.Range("D19").Value = '7.06'
$ws.Range("E19").Value = '  -2.44%  '

$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").Value = '11.03'
$ws.Range("E21").Value = '  +4.46%  '

$ws.Range("D22").Value = '469.65'
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '0.706'
$ws.Range("E23").Value = '  -2.16%  '

$ws.Range("D24").Value = '84.82'
$ws.Range("E24").Value = '  +1.23%  '

$ws.Range("E25").Value = '  -3.21%  '

$ws.Range("E26").Value = '  +0.45%  '

$ws.Range("D27").Value = '12.25'
$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("D28").Value = '10.20'
$ws.Range("E28").Value = '  -1.10%  '

$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").Value = '3.931.81'
$ws.Range("E30").Value = '  +0.53%  '

$ws.Range("E31").Value = '  -3.47%  '

$ws.Range("E32").Value = '  -3.40%  '

$ws.Range("D33").Value = '2.24'
$ws.Range("E33").Value = '  -0.81%  '

$ws.Range("D34").Value = '30.17'
$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("D35").Value = '9.36'
$ws.Range("E35").Value = '  +2.19%  '

$ws.Range("D36").Value = '1.00'

$ws.Range("D37").Value = '3.739.25'
$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("E38").Value = '  -2.97%  '

$ws.Range("E39").Value = '  -10.29%  '

$ws.Range("E40").Value = '  +1.08%  '

$ws.Range("D41").Value = '1.01'
$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").Value = '0.309'
$ws.Range("E44").Value = '  -1.30%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("E46").Value = '  +0.88%  '

$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").Value = '43.42'
$ws.Range("E47").Value = '  +11.42%  '

$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '8.61'
$ws.Range("E48").Value = '  -1.27%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '46.02'
$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value = '399.84'
$ws.Range("E50").Value = '  +0.60%  '

$ws.Range("D51").Value = '145.62'
$ws.Range("E51").Value = '  +2.53%  '
